$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 873; this shifts existing rows 873-902 down to 874-903
# and keeps their values/styles intact (matches the rest of the diff, where
# every row below is simply the previous row's data shifted down by one).
$ws.Rows.Item(873).Insert()

# Populate the newly inserted row 873 with the new data point.
$ws.Cells.Item(873, 1).Value = 4
$ws.Cells.Item(873, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(873, 3).Value = "Los Lagos"
$ws.Cells.Item(873, 4).Value = 45075
$ws.Cells.Item(873, 5).Value = 10
$ws.Cells.Item(873, 6).Value = 100112004
$ws.Cells.Item(873, 7).Value = "Cebolla"
$ws.Cells.Item(873, 8).Value = "Sin especificar"
$ws.Cells.Item(873, 9).Value = "1a (guarda)"
$ws.Cells.Item(873, 10).Value = 250
$ws.Cells.Item(873, 11).Value = 11000
$ws.Cells.Item(873, 12).Value = 11000
$ws.Cells.Item(873, 13).Value = 11000
$ws.Cells.Item(873, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(873, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(873, 16).Value = 611
$ws.Cells.Item(873, 17).Value = 18
$ws.Cells.Item(873, 18).Value = "Hortaliza"

# D873 needs the same datetime number format as the other date cells in
# column D (style index 2 / "YYYY-MM-DD HH:MM:SS"), matching the surrounding
# rows that were shifted down.
$ws.Cells.Item(873, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
